$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 40
$ws.Range("B2").Value = 33
$ws.Range("D2").Value = 7

$ws.Range("B5").Value = 0.825
$ws.Range("D5").Value = 0.175
